$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.067.42"
$ws.Range("E2").Value = "  +4.45%  "
$ws.Range("D3").Value = "3.298.00"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").Value = "630.49"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").Value = "0.413"
$ws.Range("E7").Value = "  +7.01%  "
$ws.Range("D8").Value = "0.702"
$ws.Range("E8").Value = "  +5.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "3.283.98"
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("D11").Value = "0.589"
$ws.Range("E11").Value = "  +1.59%  "
$ws.Range("D12").Value = "0.0000265"
$ws.Range("E12").Value = "  +1.57%  "
$ws.Range("D13").Value = "0.179"
$ws.Range("E13").Value = "  +1.03%  "
$ws.Range("D14").Value = "34.33"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").Value = "3.917.73"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.70%  "
$ws.Range("D17").Value = "90.881.41"
$ws.Range("E17").Value = "  +4.68%  "
$ws.Range("D18").Value = "3.312.40"
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("D19").Value = "3.22"
$ws.Range("E19").Value = "  +5.78%  "
$ws.Range("D20").Value = "14.17"
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("D21").Value = "431.17"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "8.94"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "5.34"
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("D24").Value = "0.0000177"
$ws.Range("E24").Value = "  +35.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.65%  "
$ws.Range("D26").Value = "12.16"
$ws.Range("E26").Value = "  -2.67%  "
$ws.Range("D27").Value = "3.475.30"
$ws.Range("E27").Value = "  +0.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "76.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "0.174"
$ws.Range("E30").Value = "  -2.71%  "
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("D32").Value = "559.88"
$ws.Range("E32").Value = "  +2.09%  "
$ws.Range("D33").Value = "8.65"
$ws.Range("E33").Value = "  -2.17%  "
$ws.Range("D34").Value = "7.16"
$ws.Range("E34").Value = "  +1.59%  "
$ws.Range("D35").Value = "1.35"
$ws.Range("E35").Value = "  -5.99%  "
$ws.Range("D36").Value = "1.92"
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("D37").Value = "3.56"
$ws.Range("E37").Value = "  +20.52%  "
$ws.Range("D38").Value = "22.62"
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("D39").Value = "0.134"
$ws.Range("E39").Value = "  -3.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "0.391"
$ws.Range("E42").Value = "  -1.23%  "
$ws.Range("D43").Value = "1.98"
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").Value = "184.41"
$ws.Range("E45").Value = "  +2.57%  "
$ws.Range("D46").Value = "149.08"
$ws.Range("E46").Value = "  -5.42%  "
$ws.Range("D47").Value = "44.06"
$ws.Range("E47").Value = "  -1.22%  "
$ws.Range("D48").Value = "0.129"
$ws.Range("E48").Value = "  +7.68%  "
$ws.Range("D49").Value = "1.28"
$ws.Range("E49").Value = "  -2.35%  "
$ws.Range("D50").Value = "0.628"
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("B51").Value = "Filecoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D51").Value = "4.16"
$ws.Range("E51").Value = "  -2.49%  "
